# Generate Report for Handoff
#
# b.md has been re-handed-off: mark it "Ready for handoff" on the Overview
# sheet and on each per-language sheet (zh-cn, de-de); record the new
# handoff xliff file name / timestamp and the "not latest" error detail
# that the handback/status checker produced, and widen the Error Detail
# column so the message is readable.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: row 3 is the b.md entry ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-16 02:34:10"

# ---- zh-cn sheet: row 3 is the b.md entry ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# "False" must stay a literal string, not auto-convert to a Boolean -
# prefix with an apostrophe like a user forcing text entry in Excel.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-16 02:34:02"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f39cd4c9debd6a6bea08b64d842b5393f3930a15/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d907d22b6471fcc0a370a7c88ebfe90a539133d/e2e/b.md."
# Widen the Error Detail column (P) to fit the long message - match the
# width already used by the other wide columns (G / J) in this sheet.
$wZhCn = $wsZhCn.Columns.Item(7).ColumnWidth
$wsZhCn.Columns.Item(16).ColumnWidth = $wZhCn

# ---- de-de sheet: row 3 is the b.md entry ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-16 02:34:10"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f39cd4c9debd6a6bea08b64d842b5393f3930a15/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d907d22b6471fcc0a370a7c88ebfe90a539133d/e2e/b.md."
$wDeDe = $wsDeDe.Columns.Item(7).ColumnWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $wDeDe
